# Added Test Data For Hungary/Russia/Finland Market
# Create three new worksheets (Russia, Finland, Hungary) after "Denmark",
# each cloned from the "Denmark" sheet (same layout/template) with the
# market-specific Jira reference and market name updated.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("Denmark")

# --- Russia ---------------------------------------------------------
$template.Copy($null, $template)
$russia = $wb.Worksheets.Item($wb.Worksheets.Count)
$russia.Name = "Russia"
$russia.Range("B2").Value = "NGC-2929/T2907"
$russia.Range("B4").Value = "Russia Market"
$russia.Range("A1:D11").Select()

# --- Finland ---------------------------------------------------------
$russia.Copy($null, $russia)
$finland = $wb.Worksheets.Item($wb.Worksheets.Count)
$finland.Name = "Finland"
$finland.Range("B2").Value = "NGC-3130/T2884"
$finland.Range("B4").Value = "Finland Market"
$finland.Range("A1:D11").Select()

# --- Hungary ---------------------------------------------------------
$finland.Copy($null, $finland)
$hungary = $wb.Worksheets.Item($wb.Worksheets.Count)
$hungary.Name = "Hungary"
$hungary.Range("B2").Value = "NGC-3104/T2976"
$hungary.Range("B4").Value = "Hungary Market"
$hungary.Range("H16").Select()

$hungary.Activate()
